# "add Test Case Explanation.txt update test data"
# testData.xlsx: drop the unused "actualResult" column (C), update the
# "NA" result text to "no response currently", and append a new test case
# row (input 2000 -> "no response currently").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column C ("actualResult" header, only had data in C1) entirely.
$ws.Range("C1").EntireColumn.Delete() | Out-Null

# B7 previously held "NA" -- update its text.
$ws.Range("B7").Value = "no response currently"

# Append the new test case row.
$ws.Range("A8").Value = 2000
$ws.Range("B8").Value = "no response currently"

# Match the author's final selection/active cell.
$ws.Range("B12").Select() | Out-Null
